{"js": "// The diff prefixes \"Design: \" to the start of the text of every answer\n// paragraph in the feedback table (the paragraphs styled \"List Bullet\"\n// that hold the respondent's free-text answer). There are six such\n// paragraphs in the document; each gets \"Design: \" inserted right before\n// its existing text, inside the same run (no new run/paragraph created).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style,items/text\");\nawait context.sync();\n\n// Distinctive text that begins each of the six target answer paragraphs,\n// used to find the right paragraph independent of its index, and guard\n// against accidentally double-prefixing if the script were re-applied.\nconst targets = [\n  \"We should provide a little bit more portable extinguishers.\",\n  \"Rules are fulfilled but when NB518 left,\",\n  \"Items were purchased on time.\",\n  \"First revision of the drawing was late.\",\n  \"The communication and cooperation has been really good\",\n  \"Situation was not as bad as in NB516\"\n];\n\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.style !== \"List Bullet\") {\n    continue;\n  }\n  const text = paragraph.text;\n  const matches = targets.some((t) => text.indexOf(t) === 0);\n  if (matches && text.indexOf(\"Design: \") !== 0) {\n    paragraph.insertText(\"Design: \", Word.InsertLocation.start);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The diff prefixes \"Design: \" to the start of the text of every answer\n# paragraph in the feedback table (the paragraphs styled \"List Bullet\"\n# that hold the respondent's free-text answer). There are six such\n# paragraphs in the document; each gets \"Design: \" inserted right before\n# its existing text, inside the same run (no new run/paragraph created).\n\n$d = $word.ActiveDocument\n\n# Distinctive text that begins each of the six target answer paragraphs,\n# used to find the right paragraph independent of its index, and guard\n# against accidentally double-prefixing if the script were re-applied.\n$targets = @(\n    \"We should provide a little bit more portable extinguishers.\",\n    \"Rules are fulfilled but when NB518 left,\",\n    \"Items were purchased on time.\",\n    \"First revision of the drawing was late.\",\n    \"The communication and cooperation has been really good\",\n    \"Situation was not as bad as in NB516\"\n)\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Style.NameLocal -ne \"List Bullet\") {\n        continue\n    }\n    $text = $p.Range.Text\n    foreach ($t in $targets) {\n        if ($text.StartsWith($t) -and -not $text.StartsWith(\"Design: \")) {\n            $p.Range.InsertBefore(\"Design: \")\n            break\n        }\n    }\n}\n"}
